$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string value E7420 -> E7420L (column G, rows 2:41)
$ws.Range("G2:G41").Value = "E7420L"

# Replace the =FALSE() formulas in H2:H41 with literal boolean FALSE values
$ws.Range("H2:H41").Value = $false
